$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 47.86240033333333
$ws.Range("H2").Value = 143.587201
$ws.Range("I2").Value = 0.1228118231805696
$ws.Range("J2").Value = 0.1228118231805696
$ws.Range("M2").Value = 22.78905733333334
$ws.Range("N2").Value = 68.36717200000001
$ws.Range("O2").Value = 0.387338358085063
$ws.Range("P2").Value = 0.387338358085063
$ws.Range("Q2").Value = 1090.738985307286
$ws.Range("R2").Value = 9816.650867765573
$ws.Range("S2").Value = 0.04756972994419491
$ws.Range("T2").Value = 0.04756972994419491
$ws.Range("G3").Value = 47.86240033333333
$ws.Range("H3").Value = 143.587201
$ws.Range("I3").Value = 0.1228118231805696
$ws.Range("J3").Value = 0.1228118231805696
$ws.Range("O3").Value = 0.2784908820648457
$ws.Range("P3").Value = 0.2784908820648457
$ws.Range("Q3").Value = 784.2261314435382
$ws.Range("R3").Value = 7058.035182991845
$ws.Range("S3").Value = 0.03420197296554871
$ws.Range("T3").Value = 0.0342019729655487
$ws.Range("G4").Value = 47.86240033333333
$ws.Range("H4").Value = 143.587201
$ws.Range("I4").Value = 0.1228118231805696
$ws.Range("J4").Value = 0.1228118231805696
$ws.Range("M4").Value = 11.91569833333333
$ws.Range("N4").Value = 35.747095
$ws.Range("O4").Value = 0.2025273340779806
$ws.Range("P4").Value = 0.2025273340779806
$ws.Range("Q4").Value = 570.3139238812328
$ws.Range("R4").Value = 5132.825314931095
$ws.Range("S4").Value = 0.02487275114201711
$ws.Range("T4").Value = 0.02487275114201711
$ws.Range("G5").Value = 47.86240033333333
$ws.Range("H5").Value = 143.587201
$ws.Range("I5").Value = 0.1228118231805696
$ws.Range("J5").Value = 0.1228118231805696
$ws.Range("M5").Value = 0.4016486666666667
$ws.Range("N5").Value = 1.204946
$ws.Range("O5").Value = 0.006826694619183082
$ws.Range("P5").Value = 0.006826694619183082
$ws.Range("Q5").Value = 19.22386927734956
$ws.Range("R5").Value = 173.014823496146
$ws.Range("S5").Value = 0.0008383988124788588
$ws.Range("T5").Value = 0.0008383988124788587
$ws.Range("G6").Value = 47.86240033333333
$ws.Range("H6").Value = 143.587201
$ws.Range("I6").Value = 0.1228118231805696
$ws.Range("J6").Value = 0.1228118231805696
$ws.Range("M6").Value = 7.343594
$ws.Range("N6").Value = 22.030782
$ws.Range("O6").Value = 0.1248167311529276
$ws.Range("P6").Value = 0.1248167311529276
$ws.Range("Q6").Value = 351.4820359134647
$ws.Range("R6").Value = 3163.338323221182
$ws.Range("S6").Value = 0.01532897031633004
$ws.Range("T6").Value = 0.01532897031633004
$ws.Range("I7").Value = 0.04786922362394307
$ws.Range("J7").Value = 0.04786922362394307
$ws.Range("M7").Value = 22.78905733333334
$ws.Range("N7").Value = 68.36717200000001
$ws.Range("O7").Value = 0.387338358085063
$ws.Range("P7").Value = 0.387338358085063
$ws.Range("Q7").Value = 425.144966101993
$ws.Range("R7").Value = 3826.304694917937
$ws.Range("S7").Value = 0.01854158648130482
$ws.Range("T7").Value = 0.01854158648130482
$ws.Range("I8").Value = 0.04786922362394307
$ws.Range("J8").Value = 0.04786922362394307
$ws.Range("O8").Value = 0.2784908820648457
$ws.Range("P8").Value = 0.2784908820648457
$ws.Range("S8").Value = 0.01333114231079126
$ws.Range("T8").Value = 0.01333114231079126
$ws.Range("I9").Value = 0.04786922362394307
$ws.Range("J9").Value = 0.04786922362394307
$ws.Range("M9").Value = 11.91569833333333
$ws.Range("N9").Value = 35.747095
$ws.Range("O9").Value = 0.2025273340779806
$ws.Range("P9").Value = 0.2025273340779806
$ws.Range("Q9").Value = 222.2952485444289
$ws.Range("R9").Value = 2000.65723689986
$ws.Range("S9").Value = 0.009694826244939881
$ws.Range("T9").Value = 0.009694826244939881
$ws.Range("I10").Value = 0.04786922362394307
$ws.Range("J10").Value = 0.04786922362394307
$ws.Range("M10").Value = 0.4016486666666667
$ws.Range("N10").Value = 1.204946
$ws.Range("O10").Value = 0.006826694619183082
$ws.Range("P10").Value = 0.006826694619183082
$ws.Range("Q10").Value = 7.493022035849779
$ws.Range("R10").Value = 67.437198322648
$ws.Range("S10").Value = 0.0003267885713380438
$ws.Range("T10").Value = 0.0003267885713380438
$ws.Range("I11").Value = 0.04786922362394307
$ws.Range("J11").Value = 0.04786922362394307
$ws.Range("M11").Value = 7.343594
$ws.Range("N11").Value = 22.030782
$ws.Range("O11").Value = 0.1248167311529276
$ws.Range("P11").Value = 0.1248167311529276
$ws.Range("Q11").Value = 136.9996124249573
$ws.Range("R11").Value = 1232.996511824616
$ws.Range("S11").Value = 0.005974880015569072
$ws.Range("T11").Value = 0.005974880015569073
$ws.Range("G12").Value = 171.0598806666667
$ws.Range("H12").Value = 513.1796420000001
$ws.Range("I12").Value = 0.4389285884413335
$ws.Range("J12").Value = 0.4389285884413335
$ws.Range("M12").Value = 22.78905733333334
$ws.Range("N12").Value = 68.36717200000001
$ws.Range("O12").Value = 0.387338358085063
$ws.Range("P12").Value = 0.387338358085063
$ws.Range("Q12").Value = 3898.293427945826
$ws.Range("R12").Value = 35084.64085151243
$ws.Range("S12").Value = 0.1700138787634604
$ws.Range("T12").Value = 0.1700138787634605
$ws.Range("G13").Value = 171.0598806666667
$ws.Range("H13").Value = 513.1796420000001
$ws.Range("I13").Value = 0.4389285884413335
$ws.Range("J13").Value = 0.4389285884413335
$ws.Range("O13").Value = 0.2784908820648457
$ws.Range("P13").Value = 0.2784908820648457
$ws.Range("Q13").Value = 2802.818653601583
$ws.Range("R13").Value = 25225.36788241425
$ws.Range("S13").Value = 0.1222376097585046
$ws.Range("T13").Value = 0.1222376097585046
$ws.Range("G14").Value = 171.0598806666667
$ws.Range("H14").Value = 513.1796420000001
$ws.Range("I14").Value = 0.4389285884413335
$ws.Range("J14").Value = 0.4389285884413335
$ws.Range("M14").Value = 11.91569833333333
$ws.Range("N14").Value = 35.747095
$ws.Range("O14").Value = 0.2025273340779806
$ws.Range("P14").Value = 0.2025273340779806
$ws.Range("Q14").Value = 2038.297934959999
$ws.Range("R14").Value = 18344.68141463999
$ws.Range("S14").Value = 0.0888950368676344
$ws.Range("T14").Value = 0.0888950368676344
$ws.Range("G15").Value = 171.0598806666667
$ws.Range("H15").Value = 513.1796420000001
$ws.Range("I15").Value = 0.4389285884413335
$ws.Range("J15").Value = 0.4389285884413335
$ws.Range("M15").Value = 0.4016486666666667
$ws.Range("N15").Value = 1.204946
$ws.Range("O15").Value = 0.006826694619183082
$ws.Range("P15").Value = 0.006826694619183082
$ws.Range("Q15").Value = 68.7059729899258
$ws.Range("R15").Value = 618.3537569093321
$ws.Range("S15").Value = 0.002996431432918077
$ws.Range("T15").Value = 0.002996431432918077
$ws.Range("G16").Value = 171.0598806666667
$ws.Range("H16").Value = 513.1796420000001
$ws.Range("I16").Value = 0.4389285884413335
$ws.Range("J16").Value = 0.4389285884413335
$ws.Range("M16").Value = 7.343594
$ws.Range("N16").Value = 22.030782
$ws.Range("O16").Value = 0.1248167311529276
$ws.Range("P16").Value = 0.1248167311529276
$ws.Range("Q16").Value = 1256.194313304449
$ws.Range("R16").Value = 11305.74881974005
$ws.Range("S16").Value = 0.05478563161881592
$ws.Range("T16").Value = 0.05478563161881592
$ws.Range("G17").Value = 12.628047
$ws.Range("H17").Value = 37.884141
$ws.Range("I17").Value = 0.0324027517316099
$ws.Range("J17").Value = 0.0324027517316099
$ws.Range("M17").Value = 22.78905733333334
$ws.Range("N17").Value = 68.36717200000001
$ws.Range("O17").Value = 0.387338358085063
$ws.Range("P17").Value = 0.387338358085063
$ws.Range("Q17").Value = 287.781287091028
$ws.Range("R17").Value = 2590.031583819252
$ws.Range("S17").Value = 0.01255082865315971
$ws.Range("T17").Value = 0.01255082865315971
$ws.Range("G18").Value = 12.628047
$ws.Range("H18").Value = 37.884141
$ws.Range("I18").Value = 0.0324027517316099
$ws.Range("J18").Value = 0.0324027517316099
$ws.Range("O18").Value = 0.2784908820648457
$ws.Range("P18").Value = 0.2784908820648457
$ws.Range("Q18").Value = 206.910735306356
$ws.Range("R18").Value = 1862.196617757204
$ws.Range("S18").Value = 0.009023870911064248
$ws.Range("T18").Value = 0.009023870911064248
$ws.Range("G19").Value = 12.628047
$ws.Range("H19").Value = 37.884141
$ws.Range("I19").Value = 0.0324027517316099
$ws.Range("J19").Value = 0.0324027517316099
$ws.Range("M19").Value = 11.91569833333333
$ws.Range("N19").Value = 35.747095
$ws.Range("O19").Value = 0.2025273340779806
$ws.Range("P19").Value = 0.2025273340779806
$ws.Range("Q19").Value = 150.471998591155
$ws.Range("R19").Value = 1354.247987320395
$ws.Range("S19").Value = 0.006562442924993623
$ws.Range("T19").Value = 0.006562442924993623
$ws.Range("G20").Value = 12.628047
$ws.Range("H20").Value = 37.884141
$ws.Range("I20").Value = 0.0324027517316099
$ws.Range("J20").Value = 0.0324027517316099
$ws.Range("M20").Value = 0.4016486666666667
$ws.Range("N20").Value = 1.204946
$ws.Range("O20").Value = 0.006826694619183082
$ws.Range("P20").Value = 0.006826694619183082
$ws.Range("Q20").Value = 5.072038240154001
$ws.Range("R20").Value = 45.648344161386
$ws.Range("S20").Value = 0.0002212036908929066
$ws.Range("T20").Value = 0.0002212036908929066
$ws.Range("G21").Value = 12.628047
$ws.Range("H21").Value = 37.884141
$ws.Range("I21").Value = 0.0324027517316099
$ws.Range("J21").Value = 0.0324027517316099
$ws.Range("M21").Value = 7.343594
$ws.Range("N21").Value = 22.030782
$ws.Range("O21").Value = 0.1248167311529276
$ws.Range("P21").Value = 0.1248167311529276
$ws.Range("Q21").Value = 92.73525018091802
$ws.Range("R21").Value = 834.6172516282621
$ws.Range("S21").Value = 0.004044405551499411
$ws.Range("T21").Value = 0.004044405551499412
$ws.Range("G22").Value = 139.5154473333333
$ws.Range("H22").Value = 418.546342
$ws.Range("I22").Value = 0.3579876130225438
$ws.Range("J22").Value = 0.3579876130225438
$ws.Range("M22").Value = 22.78905733333334
$ws.Range("N22").Value = 68.36717200000001
$ws.Range("O22").Value = 0.387338358085063
$ws.Range("P22").Value = 0.387338358085063
$ws.Range("Q22").Value = 3179.425528164981
$ws.Range("R22").Value = 28614.82975348483
$ws.Range("S22").Value = 0.138662334242943
$ws.Range("T22").Value = 0.1386623342429431
$ws.Range("G23").Value = 139.5154473333333
$ws.Range("H23").Value = 418.546342
$ws.Range("I23").Value = 0.3579876130225438
$ws.Range("J23").Value = 0.3579876130225438
$ws.Range("O23").Value = 0.2784908820648457
$ws.Range("P23").Value = 0.2784908820648457
$ws.Range("Q23").Value = 2285.962650783228
$ws.Range("R23").Value = 20573.66385704905
$ws.Range("S23").Value = 0.09969628611893688
$ws.Range("T23").Value = 0.09969628611893688
$ws.Range("G24").Value = 139.5154473333333
$ws.Range("H24").Value = 418.546342
$ws.Range("I24").Value = 0.3579876130225438
$ws.Range("J24").Value = 0.3579876130225438
$ws.Range("M24").Value = 11.91569833333333
$ws.Range("N24").Value = 35.747095
$ws.Range("O24").Value = 0.2025273340779806
$ws.Range("P24").Value = 0.2025273340779806
$ws.Range("Q24").Value = 1662.423983264054
$ws.Range("R24").Value = 14961.81584937649
$ws.Range("S24").Value = 0.07250227689839558
$ws.Range("T24").Value = 0.07250227689839558
$ws.Range("G25").Value = 139.5154473333333
$ws.Range("H25").Value = 418.546342
$ws.Range("I25").Value = 0.3579876130225438
$ws.Range("J25").Value = 0.3579876130225438
$ws.Range("M25").Value = 0.4016486666666667
$ws.Range("N25").Value = 1.204946
$ws.Range("O25").Value = 0.006826694619183082
$ws.Range("P25").Value = 0.006826694619183082
$ws.Range("Q25").Value = 56.03619340083689
$ws.Range("R25").Value = 504.325740607532
$ws.Range("S25").Value = 0.002443872111555195
$ws.Range("T25").Value = 0.002443872111555195
$ws.Range("G26").Value = 139.5154473333333
$ws.Range("H26").Value = 418.546342
$ws.Range("I26").Value = 0.3579876130225438
$ws.Range("J26").Value = 0.3579876130225438
$ws.Range("M26").Value = 7.343594
$ws.Range("N26").Value = 22.030782
$ws.Range("O26").Value = 0.1248167311529276
$ws.Range("P26").Value = 0.1248167311529276
$ws.Range("Q26").Value = 1024.544801944383
$ws.Range("R26").Value = 9220.903217499445
$ws.Range("S26").Value = 0.04468284365071314
$ws.Range("T26").Value = 0.04468284365071314
